$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-17: column A (index), B (name), C, D, E (in_service)
# A = index, B = name (string), C/D = numeric, E = boolean
$data = @(
    @(8,  "A8",  6,  "B8",  "line7", "C8",  14, "D8",  11, "E8",  $true),
    @(9,  "A9",  7,  "B9",  "line8", "C9",  16, "D9",  9,  "E9",  $true),
    @(10, "A10", 8,  "B10", "extr1", "C10", 5,  "D10", 12, "E10", $true),
    @(11, "A11", 9,  "B11", "extr2", "C11", 5,  "D11", 9,  "E11", $true),
    @(12, "A12", 10, "B12", "extr3", "C12", 10, "D12", 11, "E12", $false),
    @(13, "A13", 11, "B13", "extr4", "C13", 7,  "D13", 8,  "E13", $false),
    @(14, "A14", 12, "B14", "extr5", "C14", 9,  "D14", 11, "E14", $false),
    @(15, "A15", 13, "B15", "extr6", "C15", 7,  "D15", 11, "E15", $true),
    @(16, "A16", 14, "B16", "extr7", "C16", 5,  "D16", 7,  "E16", $false),
    @(17, "A17", 15, "B17", "extr8", "C17", 8,  "D17", 5,  "E17", $true)
)

# Rows 16 and 17 are brand new rows; copy the row-15 formatting down first
# so the "index" column (A) keeps the bold/centered/bordered style (s="1").
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)

foreach ($row in $data) {
    $ws.Range($row[1]).Value = $row[2]
    $ws.Range($row[3]).Value = $row[4]
    $ws.Range($row[5]).Value = $row[6]
    $ws.Range($row[7]).Value = $row[8]
    $ws.Range($row[9]).Value = $row[10]
}
